$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CS 320 Section 103 (M-W-F 2:00 - 2:50)"
$ws.Range("A1:E1").Select()
